$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-SciValue($cellRef, $literal) {
    $ws.Range($cellRef).Formula = "=" + $literal
    $ws.Range($cellRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
}

# Row 2
$ws.Range("D2").Value = 0.9424440565802336
$ws.Range("E2").Value = 0.9424440565802336

# Row 3
$ws.Range("D3").Value = 0.8973800435544202
$ws.Range("E3").Value = 0.8973800435544202

# Row 4
Set-SciValue "D4" "2.008768748498127E-08"
Set-SciValue "E4" "2.008768748498127E-08"

# Row 5
Set-SciValue "D5" "1.352348203486892E-20"
Set-SciValue "E5" "1.352348203486892E-20"

# Row 6
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = 0.4923919708852646
$ws.Range("E6").Value = 0.4923919708852646

# Row 7
$ws.Range("D7").Value = 0.9999999999154365
Set-SciValue "E7" "8.456346733964892E-11"

# Row 8
$ws.Range("D8").Value = 0.9999999999999485
Set-SciValue "E8" "5.151434834260726E-14"

# Row 9
$ws.Range("D9").Value = 0.9999999999998366
Set-SciValue "E9" "1.63424829224823E-13"

# Row 11
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = 0.9999999997253972
Set-SciValue "E11" "2.746027849553911E-10"
$ws.Range("F11").Value = 0.5809767842292786
$ws.Range("G11").Value = 0.8
